$wb = $excel.ActiveWorkbook

$wsStudies = $wb.Worksheets.Item("studies")
$wsCounts  = $wb.Worksheets.Item("counts")

# Add the new "PMID" column header to the studies sheet (col H)
$wsStudies.Range("H1").Value = "PMID"

# Add the new "notes" column header to the counts sheet (col F)
$wsCounts.Range("F1").Value = "notes"

# Update selections: studies sheet selection moves to H2 (no longer the active tab)
$wsStudies.Range("H2").Select()

# counts sheet becomes the active tab, with selection on F2
$wsCounts.Range("F2").Select()
